$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @("name","sku","category","price","cost_price","stock_quantity","unit","hsn_code","gst_rate","is_active","id","description")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Row 2: Test Widget
$ws.Cells.Item(2, 1).Value = "Test Widget"
$ws.Cells.Item(2, 2).Value = "SKU-100"
$ws.Cells.Item(2, 3).Value = "UnitTest"
$ws.Cells.Item(2, 4).Value = 99
$ws.Cells.Item(2, 5).Value = 80
$ws.Cells.Item(2, 6).Value = 10
$ws.Cells.Item(2, 7).Value = "piece"
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = "9901"
$ws.Cells.Item(2, 9).Value = 18
$ws.Cells.Item(2, 10).Value = $true
$ws.Cells.Item(2, 11).Value = "5d3821c0-8d44-4733-8fe3-786c898948b1"
$ws.Cells.Item(2, 12).NumberFormat = "@"
$ws.Cells.Item(2, 12).Value = ""

# Row 3: Hassan mansuri
$ws.Cells.Item(3, 1).Value = "Hassan mansuri"
$ws.Cells.Item(3, 2).NumberFormat = "@"
$ws.Cells.Item(3, 2).Value = "4"
$ws.Cells.Item(3, 3).NumberFormat = "@"
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = "piece"
$ws.Cells.Item(3, 8).NumberFormat = "@"
$ws.Cells.Item(3, 8).Value = "3243"
$ws.Cells.Item(3, 9).Value = 18
$ws.Cells.Item(3, 10).Value = $true
$ws.Cells.Item(3, 11).Value = "b8a778da-786f-4741-80a3-f42a37f8bde6"
$ws.Cells.Item(3, 12).Value = "dfgrhah regerge"

# Row 4: Test Widget (duplicate, no description cell)
$ws.Cells.Item(4, 1).Value = "Test Widget"
$ws.Cells.Item(4, 2).Value = "SKU-100"
$ws.Cells.Item(4, 3).Value = "UnitTest"
$ws.Cells.Item(4, 4).Value = 99
$ws.Cells.Item(4, 5).Value = 80
$ws.Cells.Item(4, 6).Value = 10
$ws.Cells.Item(4, 7).Value = "piece"
$ws.Cells.Item(4, 8).NumberFormat = "@"
$ws.Cells.Item(4, 8).Value = "9901"
$ws.Cells.Item(4, 9).Value = 18
$ws.Cells.Item(4, 10).Value = $true
$ws.Cells.Item(4, 11).Value = "3f441215-d35d-4111-acd6-cdcf5fa68555"
